$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Yohannes Fantahun", "fantish@mail.com", "M", "A", "-", 2, "Unity University College"),
    @(5, "test person", "test@mail.com", "M", "A", "-", 0, "Addis Ababa University"),
    @(6, "Abex Abelew", "abex@mail.com", "M", "e", "-", 1, "Addis Ababa University")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
